$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.449.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "'1.567.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D6").Value = "'288.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'0.3720"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").Value = "'48.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").Value = "'0.3318"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "'0.07469"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "'1.128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "'5.956"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'6.902"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "'1.579.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06765"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'87.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'6.351"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'16.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'12.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").Value = "'22.437.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'2.393"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").Value = "'2.556"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("D27").Value = "'153.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "'19.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").Value = "'5.019"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "'123.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "'1.744.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'1.054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "'2.014"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'6.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'9.605"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "'0.08299"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'0.02452"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "'0.2268"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "'0.06378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").Value = "'5.366"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'1.289"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("D42").Value = "'0.6272"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "'1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "'13.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "'0.6122"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("D47").Value = "'3.783"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'2.041"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'125.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "'1.213"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("D51").Value = "'0.07240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
